$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 66: copy formatting from row 65, then set values ---
$ws.Range("A65:V65").Copy() | Out-Null
$ws.Range("A66").PasteSpecial(-4122) | Out-Null

# Row 66 values
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = 'scotland'
$ws.Range("C66").Value = 'league-one'
$ws.Range("D66").Value = '2023-2024'
$ws.Range("E66").Value = 45244.86458333334
$ws.Range("F66").Value = 'Cove Rangers'
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 'Montrose'
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 1.83
$ws.Range("K66").Value = '19/10/2023 09:13'
$ws.Range("L66").Value = 1.73
$ws.Range("M66").Value = '14/11/2023 20:36'
$ws.Range("N66").Value = 3.63
$ws.Range("O66").Value = '19/10/2023 09:13'
$ws.Range("P66").Value = 3.87
$ws.Range("Q66").Value = '14/11/2023 20:36'
$ws.Range("R66").Value = 3.59
$ws.Range("S66").Value = '19/10/2023 09:13'
$ws.Range("T66").Value = 4.51
$ws.Range("U66").Value = '14/11/2023 20:36'
$ws.Range("V66").Value = 'https://www.betexplorer.com/football/scotland/league-one/cove-rangers-montrose/feg2EOWC/'

# Row 2 (F:V updated)
$ws.Range("F2").Value = 'Montrose'
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 'Kelty Hearts'
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 1.95
$ws.Range("K2").Value = '04/08/2023 04:12'
$ws.Range("L2").Value = 2.04
$ws.Range("M2").Value = '05/08/2023 15:58'
$ws.Range("N2").Value = 3.35
$ws.Range("O2").Value = '04/08/2023 04:12'
$ws.Range("P2").Value = 3.8
$ws.Range("Q2").Value = '05/08/2023 15:58'
$ws.Range("R2").Value = 3.44
$ws.Range("S2").Value = '04/08/2023 04:12'
$ws.Range("T2").Value = 3.29
$ws.Range("U2").Value = '05/08/2023 15:58'
$ws.Range("V2").Value = 'https://www.betexplorer.com/football/scotland/league-one/montrose-kelty-hearts/0n7822Mt/'

# Row 5 (F:V updated)
$ws.Range("F5").Value = 'Hamilton'
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 'Cove Rangers'
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1.83
$ws.Range("K5").Value = '04/08/2023 16:33'
$ws.Range("L5").Value = 1.9
$ws.Range("M5").Value = '05/08/2023 15:57'
$ws.Range("N5").Value = 3.77
$ws.Range("O5").Value = '04/08/2023 16:33'
$ws.Range("P5").Value = 3.76
$ws.Range("Q5").Value = '05/08/2023 15:57'
$ws.Range("R5").Value = 3.95
$ws.Range("S5").Value = '04/08/2023 16:33'
$ws.Range("T5").Value = 3.75
$ws.Range("U5").Value = '05/08/2023 15:57'
$ws.Range("V5").Value = 'https://www.betexplorer.com/football/scotland/league-one/hamilton-cove-rangers/bwLt7rEP/'

# Row 47 (F:V updated)
$ws.Range("F47").Value = 'Stirling'
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 'Falkirk'
$ws.Range("I47").Value = 2
$ws.Range("J47").Value = 5.43
$ws.Range("K47").Value = '19/10/2023 08:13'
$ws.Range("L47").Value = 6.82
$ws.Range("M47").Value = '21/10/2023 15:56'
$ws.Range("N47").Value = 4.14
$ws.Range("O47").Value = '19/10/2023 08:13'
$ws.Range("P47").Value = 4.59
$ws.Range("Q47").Value = '21/10/2023 15:56'
$ws.Range("R47").Value = 1.5
$ws.Range("S47").Value = '19/10/2023 08:13'
$ws.Range("T47").Value = 1.44
$ws.Range("U47").Value = '21/10/2023 15:56'
$ws.Range("V47").Value = 'https://www.betexplorer.com/football/scotland/league-one/stirling-falkirk/ziWlxSOJ/'

# Row 48 (F:V updated)
$ws.Range("F48").Value = 'Hamilton'
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 'Edinburgh City'
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 1.19
$ws.Range("K48").Value = '19/10/2023 08:13'
$ws.Range("L48").Value = 1.17
$ws.Range("M48").Value = '21/10/2023 15:35'
$ws.Range("N48").Value = 6.55
$ws.Range("O48").Value = '19/10/2023 08:13'
$ws.Range("P48").Value = 7.9
$ws.Range("Q48").Value = '21/10/2023 15:35'
$ws.Range("R48").Value = 9.529999999999999
$ws.Range("S48").Value = '19/10/2023 08:13'
$ws.Range("T48").Value = 13.65
$ws.Range("U48").Value = '21/10/2023 15:35'
$ws.Range("V48").Value = 'https://www.betexplorer.com/football/scotland/league-one/hamilton-edinburgh-city/Yoh6DrnJ/'

# Row 57 (F:V updated)
$ws.Range("F57").Value = 'Queen of South'
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 'Falkirk'
$ws.Range("I57").Value = 1
$ws.Range("J57").Value = 5.39
$ws.Range("K57").Value = '02/11/2023 08:13'
$ws.Range("L57").Value = 8.18
$ws.Range("M57").Value = '04/11/2023 15:36'
$ws.Range("N57").Value = 4.14
$ws.Range("O57").Value = '02/11/2023 08:13'
$ws.Range("P57").Value = 5.09
$ws.Range("Q57").Value = '04/11/2023 15:36'
$ws.Range("R57").Value = 1.5
$ws.Range("S57").Value = '02/11/2023 08:13'
$ws.Range("T57").Value = 1.36
$ws.Range("U57").Value = '04/11/2023 15:36'
$ws.Range("V57").Value = 'https://www.betexplorer.com/football/scotland/league-one/queen-of-south-falkirk/hlpD9Zpo/'

# Row 59 (F:V updated)
$ws.Range("F59").Value = 'Annan'
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 'Kelty Hearts'
$ws.Range("I59").Value = 2
$ws.Range("J59").Value = 2.5
$ws.Range("K59").Value = '02/11/2023 08:13'
$ws.Range("L59").Value = 2.45
$ws.Range("M59").Value = '04/11/2023 15:54'
$ws.Range("N59").Value = 3.33
$ws.Range("O59").Value = '02/11/2023 08:13'
$ws.Range("P59").Value = 3.69
$ws.Range("Q59").Value = '04/11/2023 15:54'
$ws.Range("R59").Value = 2.5
$ws.Range("S59").Value = '02/11/2023 08:13'
$ws.Range("T59").Value = 2.64
$ws.Range("U59").Value = '04/11/2023 15:54'
$ws.Range("V59").Value = 'https://www.betexplorer.com/football/scotland/league-one/annan-kelty-hearts/zFEyTNo6/'

# Row 60 (F:V updated)
$ws.Range("F60").Value = 'Edinburgh City'
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 'Montrose'
$ws.Range("I60").Value = 5
$ws.Range("J60").Value = 3.27
$ws.Range("K60").Value = '02/11/2023 08:13'
$ws.Range("L60").Value = 3.05
$ws.Range("M60").Value = '04/11/2023 14:22'
$ws.Range("N60").Value = 3.59
$ws.Range("O60").Value = '02/11/2023 08:13'
$ws.Range("P60").Value = 3.85
$ws.Range("Q60").Value = '04/11/2023 15:34'
$ws.Range("R60").Value = 1.93
$ws.Range("S60").Value = '02/11/2023 08:13'
$ws.Range("T60").Value = 2.13
$ws.Range("U60").Value = '04/11/2023 15:34'
$ws.Range("V60").Value = 'https://www.betexplorer.com/football/scotland/league-one/edinburgh-city-montrose/GOMqR1GI/'

# Row 61 (F:V updated)
$ws.Range("F61").Value = 'Montrose'
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 'Annan'
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 1.9
$ws.Range("K61").Value = '09/11/2023 09:12'
$ws.Range("L61").Value = 1.79
$ws.Range("M61").Value = '11/11/2023 15:54'
$ws.Range("N61").Value = 3.92
$ws.Range("O61").Value = '09/11/2023 09:12'
$ws.Range("P61").Value = 4.32
$ws.Range("Q61").Value = '11/11/2023 15:54'
$ws.Range("R61").Value = 3.22
$ws.Range("S61").Value = '09/11/2023 09:12'
$ws.Range("T61").Value = 3.74
$ws.Range("U61").Value = '11/11/2023 15:54'
$ws.Range("V61").Value = 'https://www.betexplorer.com/football/scotland/league-one/montrose-annan/lzUFTfxH/'

# Row 62 (F:V updated)
$ws.Range("F62").Value = 'Alloa'
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 'Queen of South'
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2.41
$ws.Range("K62").Value = '09/11/2023 09:12'
$ws.Range("L62").Value = 2.46
$ws.Range("M62").Value = '11/11/2023 15:49'
$ws.Range("N62").Value = 3.29
$ws.Range("O62").Value = '09/11/2023 09:12'
$ws.Range("P62").Value = 3.52
$ws.Range("Q62").Value = '11/11/2023 15:20'
$ws.Range("R62").Value = 2.62
$ws.Range("S62").Value = '09/11/2023 09:12'
$ws.Range("T62").Value = 2.72
$ws.Range("U62").Value = '11/11/2023 15:49'
$ws.Range("V62").Value = 'https://www.betexplorer.com/football/scotland/league-one/alloa-queen-of-south/KfxaXxwh/'

# Row 63 (F:V updated)
$ws.Range("F63").Value = 'Falkirk'
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 'Edinburgh City'
$ws.Range("I63").Value = 1
$ws.Range("J63").Value = 1.14
$ws.Range("K63").Value = '09/11/2023 09:12'
$ws.Range("L63").Value = 1.09
$ws.Range("M63").Value = '10/11/2023 16:03'
$ws.Range("N63").Value = 7.57
$ws.Range("O63").Value = '09/11/2023 09:12'
$ws.Range("P63").Value = 11.04
$ws.Range("Q63").Value = '11/11/2023 15:46'
$ws.Range("R63").Value = 11.46
$ws.Range("S63").Value = '09/11/2023 09:12'
$ws.Range("T63").Value = 22.57
$ws.Range("U63").Value = '11/11/2023 15:46'
$ws.Range("V63").Value = 'https://www.betexplorer.com/football/scotland/league-one/falkirk-edinburgh-city/2JR3Wdhb/'

# Row 64 (F:V updated)
$ws.Range("F64").Value = 'Hamilton'
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 'Stirling'
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1.37
$ws.Range("K64").Value = '09/11/2023 09:12'
$ws.Range("L64").Value = 1.37
$ws.Range("M64").Value = '11/11/2023 15:52'
$ws.Range("N64").Value = 4.59
$ws.Range("O64").Value = '09/11/2023 09:12'
$ws.Range("P64").Value = 4.79
$ws.Range("Q64").Value = '11/11/2023 15:52'
$ws.Range("R64").Value = 6.39
$ws.Range("S64").Value = '09/11/2023 09:12'
$ws.Range("T64").Value = 8.41
$ws.Range("U64").Value = '11/11/2023 15:52'
$ws.Range("V64").Value = 'https://www.betexplorer.com/football/scotland/league-one/hamilton-stirling/GMV7VG74/'

# Row 65 (F:V updated)
$ws.Range("F65").Value = 'Kelty Hearts'
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 'Cove Rangers'
$ws.Range("I65").Value = 1
$ws.Range("J65").Value = 2.46
$ws.Range("K65").Value = '09/11/2023 09:12'
$ws.Range("L65").Value = 3.06
$ws.Range("M65").Value = '11/11/2023 15:57'
$ws.Range("N65").Value = 3.33
$ws.Range("O65").Value = '09/11/2023 09:12'
$ws.Range("P65").Value = 3.51
$ws.Range("Q65").Value = '11/11/2023 15:57'
$ws.Range("R65").Value = 2.53
$ws.Range("S65").Value = '09/11/2023 09:12'
$ws.Range("T65").Value = 2.24
$ws.Range("U65").Value = '11/11/2023 15:53'
$ws.Range("V65").Value = 'https://www.betexplorer.com/football/scotland/league-one/kelty-hearts-cove-rangers/AqVBUzNA/'
